# Apply updated crypto price / volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; D = '61.901.63'; E = '  -0.70%  ' }
    @{ Row = 3; D = '2.409.11'; E = '  -0.62%  ' }
    @{ Row = 4; D = $null; E = '  -0.02%  ' }
    @{ Row = 5; D = '561.96'; E = '  +1.03%  ' }
    @{ Row = 6; D = '142.17'; E = '  -1.10%  ' }
    @{ Row = 7; D = $null; E = '  +0.02%  ' }
    @{ Row = 8; D = $null; E = '  -0.90%  ' }
    @{ Row = 9; D = $null; E = '  -0.81%  ' }
    @{ Row = 10; D = $null; E = '  -1.82%  ' }
    @{ Row = 11; D = $null; E = '  -2.93%  ' }
    @{ Row = 12; D = '0.348'; E = '  -0.81%  ' }
    @{ Row = 13; D = '25.48'; E = '  -2.97%  ' }
    @{ Row = 14; D = $null; E = '  -1.62%  ' }
    @{ Row = 15; D = '2.843.03'; E = '  -0.65%  ' }
    @{ Row = 16; D = '61.895.51'; E = '  -0.51%  ' }
    @{ Row = 17; D = '2.405.64'; E = '  -0.72%  ' }
    @{ Row = 18; D = $null; E = '  +1.25%  ' }
    @{ Row = 19; D = '320.80'; E = '  -1.12%  ' }
    @{ Row = 20; D = '6.82'; E = '  +1.25%  ' }
    @{ Row = 21; D = $null; E = '  -1.75%  ' }
    @{ Row = 22; D = $null; E = '  -0.18%  ' }
    @{ Row = 23; D = '65.51'; E = '  +1.07%  ' }
    @{ Row = 24; D = $null; E = '  -2.83%  ' }
    @{ Row = 25; D = $null; E = '  -4.84%  ' }
    @{ Row = 26; D = '564.54'; E = '  -1.46%  ' }
    @{ Row = 27; D = '1.00'; E = '  -0.03%  ' }
    @{ Row = 28; D = '2.515.36'; E = '  -0.92%  ' }
    @{ Row = 29; D = '0.0₃0931'; E = '  -0.87%  ' }
    @{ Row = 30; D = $null; E = '  -2.72%  ' }
    @{ Row = 31; D = $null; E = '  -4.84%  ' }
    @{ Row = 32; D = $null; E = '  -0.83%  ' }
    @{ Row = 33; D = '1.87'; E = '  +0.33%  ' }
    @{ Row = 34; D = $null; E = '  -4.09%  ' }
    @{ Row = 35; D = $null; E = '  +0.02%  ' }
    @{ Row = 36; D = '4.74'; E = '  -1.76%  ' }
    @{ Row = 37; D = $null; E = '  -4.94%  ' }
    @{ Row = 38; D = '152.45'; E = '  +1.89%  ' }
    @{ Row = 39; D = $null; E = '  -1.51%  ' }
    @{ Row = 40; D = '18.48'; E = '  -1.52%  ' }
    @{ Row = 41; D = $null; E = '  -5.13%  ' }
    @{ Row = 42; D = $null; E = '  -0.03%  ' }
    @{ Row = 43; D = '147.91'; E = '  -2.02%  ' }
    @{ Row = 44; D = $null; E = '  -3.90%  ' }
    @{ Row = 45; D = $null; E = '  -1.31%  ' }
    @{ Row = 46; D = '0.0528'; E = '  -2.74%  ' }
    @{ Row = 47; D = '19.81'; E = '  -2.67%  ' }
    @{ Row = 48; D = $null; E = '  +0.01%  ' }
    @{ Row = 49; D = '0.0918'; E = '  +0.25%  ' }
    @{ Row = 50; D = $null; E = '  -1.85%  ' }
)

foreach ($change in $changes) {
    $row = $change.Row
    if ($null -ne $change.D) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $change.D
        $cellD.Style = "Normal"
    }
    if ($null -ne $change.E) {
        $ws.Cells.Item($row, 5).Value = $change.E
    }
}
